$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order")

# Insert a new row before row 482 (shifts existing row 482 and below down by one)
$ws.Rows.Item(482).Insert()

# Populate the newly inserted row 482 with the new field definition
$ws.Cells.Item(482, 1).Value = "Order"
$ws.Cells.Item(482, 2).Value = "Site_Name_No_Spaces__c"
$ws.Cells.Item(482, 3).Value = "string"
$ws.Cells.Item(482, 4).Value = "Site Name No Spaces"
$ws.Cells.Item(482, 5).Value = 1300
$ws.Cells.Item(482, 6).Value = $true
$ws.Cells.Item(482, 7).Value = $true
$ws.Cells.Item(482, 8).Value = $true
